# Added support for Java
# - Header row: questions were re-labelled / reordered (Ques3,Ques1,Ques2 -> Ques1,Ques2,Ques3)
# - Row 2 (Darshan_Padia_65): Total Marks becomes a textual "12.0" instead of numeric 22
# - Row 3 (Mustafa_Bharamal_78): Ques1/Ques2 percentages swapped, Total Marks becomes textual "5.7"
# - Row 4 (Priya_Rajani_12): Ques3 becomes "100.0%", Total Marks becomes textual "4.0"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: rotate labels left (B<-Ques1, C<-Ques2, D<-Ques3) ---
$ws.Range("B1").Value = "Ques1"
$ws.Range("C1").Value = "Ques2"
$ws.Range("D1").Value = "Ques3"

# --- Row 2 ---
# Force the value to be stored as literal text (even though it looks like a
# number) by temporarily marking the cell as Text, assigning it, then
# clearing the formatting again so no stray style is left behind.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "12.0"
$ws.Range("E2").ClearFormats()

# --- Row 3: swap Ques1 (B3) and Ques2 (C3) percentages ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "66.7%"
$ws.Range("B3").ClearFormats()

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "33.3%"
$ws.Range("C3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.7"
$ws.Range("E3").ClearFormats()

# --- Row 4 ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "100.0%"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.0"
$ws.Range("E4").ClearFormats()
